$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New block 1: rows 18-24, duplicate of the 9-15 per-fold table, numeric values ---

# Row 18 header (reuse existing shared strings via format-copy from row 9's header style)
$ws.Range("B9:F9").Copy()
$ws.Range("B18:F18").PasteSpecial(-4122)
$ws.Range("B18").Value = "name_recall"
$ws.Range("C18").Value = "name_precision"
$ws.Range("D18").Value = "type_accuracy"
$ws.Range("E18").Value = "value_recall"
$ws.Range("F18").Value = "value_precision"

# Rows 19-24: same per-row label + numeric values as rows 10-15, formatted 0.0000_
$ws.Range("B19:F24").NumberFormat = "0.0000_ "

$ws.Range("A19").Value = "nodes_0shot_0"
$ws.Range("B19").Value = 0.73373401140391403
$ws.Range("C19").Value = 0.79257589767298475
$ws.Range("D19").Value = 0.71558021266759109
$ws.Range("E19").Value = 0.86812297734627841
$ws.Range("F19").Value = 0.89158576051779947

$ws.Range("A20").Value = "nodes_1shot_0"
$ws.Range("B20").Value = 0.77297349360456169
$ws.Range("C20").Value = 0.81281399291108025
$ws.Range("D20").Value = 0.6920711974110032
$ws.Range("E20").Value = 0.91100323624595481
$ws.Range("F20").Value = 0.92071197411003247

$ws.Range("A21").Value = "nodes_3shot_0"
$ws.Range("B21").Value = 0.83350670365233437
$ws.Range("C21").Value = 0.76533329130416505
$ws.Range("D21").Value = 0.67102018801047913
$ws.Range("E21").Value = 0.89482200647249199
$ws.Range("F21").Value = 0.89482200647249199

$ws.Range("A22").Value = "nodes_3shot_1"
$ws.Range("B22").Value = 0.84985745107104316
$ws.Range("C22").Value = 0.74554391132061015
$ws.Range("D22").Value = 0.74323470488519039
$ws.Range("E22").Value = 0.92233009708737879
$ws.Range("F22").Value = 0.91747572815533984

$ws.Range("A23").Value = "rule_cot"
$ws.Range("B23").Value = 0.80248112189859744
$ws.Range("C23").Value = 0.7348120595693407
$ws.Range("D23").Value = 0.66680536292186776
$ws.Range("E23").Value = 0.85760517799352742
$ws.Range("F23").Value = 0.88025889967637549

$ws.Range("A24").Value = "zs_cot"
$ws.Range("B24").Value = 0.59940668824163978
$ws.Range("C24").Value = 0.79729542302357803
$ws.Range("D24").Value = 0.5553513638465094
$ws.Range("E24").Value = 0.85760517799352742
$ws.Range("F24").Value = 0.87378640776699024

# --- New block 2: rows 28-34, summary table of mean +/- std across folds, as text ---

$ws.Range("A28").NumberFormat = "@"
$ws.Range("A29:F34").NumberFormat = "@"

$ws.Range("B9:F9").Copy()
$ws.Range("B28:F28").PasteSpecial(-4122)
$ws.Range("B28:F28").NumberFormat = "@"
$ws.Range("B28").Value = "name_recall"
$ws.Range("C28").Value = "name_precision"
$ws.Range("D28").Value = "type_accuracy"
$ws.Range("E28").Value = "value_recall"
$ws.Range("F28").Value = "value_precision"

$ws.Range("A29").Value = "nodes_0shot_0"
$ws.Range("B29").Value = "0.7337±0.0391"
$ws.Range("C29").Value = "0.7926±0.0412"
$ws.Range("D29").Value = "0.7156±0.0558"
$ws.Range("E29").Value = "0.8681±0.0452"
$ws.Range("F29").Value = "0.8916±0.0427"

$ws.Range("A30").Value = "nodes_1shot_0"
$ws.Range("B30").Value = "0.7730±0.0378"
$ws.Range("C30").Value = "0.8128±0.0403"
$ws.Range("D30").Value = "0.6921±0.0572"
$ws.Range("E30").Value = "0.9110±0.0405"
$ws.Range("F30").Value = "0.9207±0.0374"

$ws.Range("A31").Value = "nodes_3shot_0"
$ws.Range("B31").Value = "0.8335±0.0371"
$ws.Range("C31").Value = "0.7653±0.0445"
$ws.Range("D31").Value = "0.6710±0.0491"
$ws.Range("E31").Value = "0.8948±0.0413"
$ws.Range("F31").Value = "0.8948±0.0413"

$ws.Range("A32").Value = "nodes_3shot_1"
$ws.Range("B32").Value = "0.8499±0.0356"
$ws.Range("C32").Value = "0.7455±0.0422"
$ws.Range("D32").Value = "0.7432±0.0538"
$ws.Range("E32").Value = "0.9223±0.0385"
$ws.Range("F32").Value = "0.9175±0.0395"

$ws.Range("A33").Value = "rule_cot"
$ws.Range("B33").Value = "0.8025±0.0462"
$ws.Range("C33").Value = "0.7348±0.0486"
$ws.Range("D33").Value = "0.6668±0.0588"
$ws.Range("E33").Value = "0.8576±0.0558"
$ws.Range("F33").Value = "0.8803±0.0545"

$ws.Range("A34").Value = "zs_cot"
$ws.Range("B34").Value = "0.5994±0.0540"
$ws.Range("C34").Value = "0.7973±0.0501"
$ws.Range("D34").Value = "0.5554±0.0760"
$ws.Range("E34").Value = "0.8576±0.0549"
$ws.Range("F34").Value = "0.8738±0.0536"

# --- View updates: selection on J18 (topLeftCell scroll position not persisted by this host) ---
$ws.Range("J18").Select() | Out-Null
